# Crypto price/volume refresh from the Tue Oct 24 15:16:43 UTC 2023 GitHub Actions run.
# Column D ("Price") values are free-form text (e.g. "34.311.69", "1.00", "0.0₃0758") that
# happen to look numeric, so we force them through as text (NumberFormat "@") and then reset
# the cell style back to Normal so no stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "34.311.69"
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +11.26%  "

$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.820.44"
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +8.37%  "

$dCell = $ws.Cells.Item(4, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.15%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "229.70"
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +4.77%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.577"
$dCell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +8.99%  "

$dCell = $ws.Cells.Item(7, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.09%  "

$dCell = $ws.Cells.Item(8, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "32.27"
$dCell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +10.50%  "

$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "46.78"
$dCell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +6.04%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.288"
$dCell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +9.39%  "

$dCell = $ws.Cells.Item(11, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0678"
$dCell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +5.17%  "

$ws.Cells.Item(12, 5).Value = "  +3.01%  "

$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.081.74"
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +8.46%  "

$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.818.38"
$dCell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +8.80%  "

$ws.Cells.Item(15, 5).Value = "  +7.24%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "34.317.57"
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +11.44%  "

$dCell = $ws.Cells.Item(17, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "10.26"
$dCell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.70%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "4.28"
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +5.98%  "

$dCell = $ws.Cells.Item(19, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "70.54"
$dCell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +6.93%  "

$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "259.11"
$dCell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +6.34%  "

$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0758"
$dCell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +5.18%  "

$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "10.69"
$dCell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +7.05%  "

$dCell = $ws.Cells.Item(24, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "4.37"
$dCell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +3.01%  "

$ws.Cells.Item(25, 5).Value = "  +2.93%  "

$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "159.95"
$dCell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.41%  "

$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "16.82"
$dCell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +6.33%  "

$ws.Cells.Item(28, 5).Value = "  +5.09%  "

$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "7.09"
$dCell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +5.96%  "

$dCell = $ws.Cells.Item(30, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.09%  "

$dCell = $ws.Cells.Item(31, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.93"
$dCell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +13.33%  "

$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0530"
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +7.32%  "

$ws.Cells.Item(33, 5).Value = "  +6.55%  "

$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.62"
$dCell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +9.72%  "

$dCell = $ws.Cells.Item(35, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.526.70"
$dCell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.04%  "

$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.81"
$dCell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +2.55%  "

$ws.Cells.Item(37, 5).Value = "  +5.27%  "

$dCell = $ws.Cells.Item(38, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.639"
$dCell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +6.13%  "

$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0191"
$dCell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +6.80%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "83.55"
$dCell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.74%  "

$ws.Cells.Item(41, 5).Value = "  +5.00%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.37"
$dCell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +3.46%  "

$dCell = $ws.Cells.Item(43, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.914"
$dCell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +9.14%  "

$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.13"
$dCell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +4.90%  "

$ws.Cells.Item(45, 5).Value = "  +4.50%  "

$ws.Cells.Item(46, 5).Value = "  +5.44%  "

$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.976.41"
$dCell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +9.01%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "5.87"
$dCell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +6.12%  "

$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "12.01"
$dCell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +14.84%  "

$ws.Cells.Item(50, 5).Value = "  -0.05%  "

$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "52.02"
$dCell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +3.14%  "
